$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 97.25
$ws.Range("I4").Value = 97.25
$ws.Range("K4").Value = 97.25
$ws.Range("M4").Value = 16.75
$ws.Range("H33").Value = 283.39285
$ws.Range("I33").Value = 296.1905
$ws.Range("J33").Value = 245
$ws.Range("K33").Value = 296.1905
$ws.Range("L33").Value = 245
$ws.Range("M33").Value = -67.19049999999999
$ws.Range("N33").Value = -703
$ws.Range("H40").Value = 2999.875
$ws.Range("I40").Value = 3425
$ws.Range("J40").Value = 2574.75
$ws.Range("K40").Value = 3425
$ws.Range("L40").Value = 2574.75
$ws.Range("M40").Value = -3250
$ws.Range("N40").Value = -2924.75
$ws.Range("H64").Value = 8740
$ws.Range("I64").Value = 3880
$ws.Range("J64").Value = 13600
$ws.Range("K64").Value = 3880
$ws.Range("L64").Value = 13600
$ws.Range("M64").Value = -3632
$ws.Range("N64").Value = -14096
$ws.Range("H67").Value = 8740
$ws.Range("I67").Value = 3880
$ws.Range("J67").Value = 13600
$ws.Range("K67").Value = 3880
$ws.Range("L67").Value = 13600
$ws.Range("M67").Value = -3022
$ws.Range("N67").Value = -15316
$ws.Range("H97").Value = 32088
$ws.Range("J97").Value = 44673.2
$ws.Range("L97").Value = 134019.6
$ws.Range("N97").Value = -135011.6
$ws.Range("H107").Value = 656.25
$ws.Range("I107").Value = 616.7778
$ws.Range("K107").Value = 616.7778
$ws.Range("M107").Value = 1303.2222
$ws.Range("H135").Value = 1226.04
$ws.Range("I135").Value = 872.5
$ws.Range("J135").Value = 2135.1428
$ws.Range("K135").Value = 7852.5
$ws.Range("L135").Value = 19216.2852
$ws.Range("M135").Value = -5317.5
$ws.Range("N135").Value = -24286.2852

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7664.082
$ws.Range("I32").Value = 3789.46
$ws.Range("K32").Value = 3789.46
$ws.Range("M32").Value = -3502.46
$ws.Range("H101").Value = 30602
$ws.Range("J101").Value = 30602
$ws.Range("L101").Value = 30602
$ws.Range("N101").Value = -37092
$ws.Range("H102").Value = 2505.0386
$ws.Range("I102").Value = 2012.6666
$ws.Range("J102").Value = 3612.875
$ws.Range("K102").Value = 2012.6666
$ws.Range("L102").Value = 3612.875
$ws.Range("M102").Value = -390.6666
$ws.Range("N102").Value = -6856.875
$ws.Range("H110").Value = 12341.944
$ws.Range("I110").Value = 13410.333
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 13410.333
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = -11365.333
$ws.Range("N110").Value = -11090

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 8574.416999999999
$ws.Range("I29").Value = 10047.3
$ws.Range("J29").Value = 1210
$ws.Range("K29").Value = 10047.3
$ws.Range("L29").Value = 1210
$ws.Range("M29").Value = -9758.299999999999
$ws.Range("N29").Value = -1788
$ws.Range("H36").Value = 3474.75
$ws.Range("I36").Value = 1633
$ws.Range("K36").Value = 1633
$ws.Range("M36").Value = -1099
$ws.Range("H75").Value = 26666.334
$ws.Range("I75").Value = 14999.5
$ws.Range("K75").Value = 14999.5
$ws.Range("M75").Value = -14063.5
$ws.Range("H78").Value = 26666.334
$ws.Range("I78").Value = 14999.5
$ws.Range("K78").Value = 44998.5
$ws.Range("M78").Value = -40318.5
$ws.Range("H86").Value = 2997.1667
$ws.Range("I86").Value = 2809.3125
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 2809.3125
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -1686.3125
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 2997.1667
$ws.Range("I89").Value = 2809.3125
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 14046.5625
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -8430.5625
$ws.Range("N89").Value = -33732
$ws.Range("H105").Value = 2053.5
$ws.Range("I105").Value = 2053.5
$ws.Range("K105").Value = 2053.5
$ws.Range("M105").Value = -306.5
$ws.Range("H134").Value = 1871.1666
$ws.Range("I134").Value = 1861.7715
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 5585.3145
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -3050.3145
$ws.Range("N134").Value = -11670
$ws.Range("H135").Value = 90300
$ws.Range("J135").Value = 90300
$ws.Range("L135").Value = 90300
$ws.Range("N135").Value = -100440

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I58").Value = 3432.9092
$ws.Range("J58").Value = 3980.4
$ws.Range("K58").Value = 3432.9092
$ws.Range("L58").Value = 3980.4
$ws.Range("M58").Value = -3229.9092
$ws.Range("N58").Value = -4386.4
$ws.Range("H60").Value = 17128.428
$ws.Range("I60").Value = 14333
$ws.Range("K60").Value = 14333
$ws.Range("M60").Value = -13822
$ws.Range("H105").Value = 1329
$ws.Range("I105").Value = 1362
$ws.Range("J105").Value = 1308.375
$ws.Range("K105").Value = 1362
$ws.Range("L105").Value = 1308.375
$ws.Range("M105").Value = 385
$ws.Range("N105").Value = -4802.375
$ws.Range("I136").Value = 3432.9092
$ws.Range("J136").Value = 3980.4
$ws.Range("K136").Value = 10298.7276
$ws.Range("L136").Value = 11941.2
$ws.Range("M136").Value = -7748.7276
$ws.Range("N136").Value = -17041.2

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 341.86957
$ws.Range("J23").Value = 354.23077
$ws.Range("L23").Value = 1062.69231
$ws.Range("N23").Value = -1532.69231
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H118").Value = 3000
$ws.Range("I118").Value = 3000
$ws.Range("K118").Value = 9000
$ws.Range("M118").Value = -7757
$ws.Range("H119").Value = 1000
$ws.Range("I119").Value = 1000
$ws.Range("K119").Value = 3000
$ws.Range("M119").Value = 1838

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 42340.668
$ws.Range("J34").Value = 42340.668
$ws.Range("L34").Value = 42340.668
$ws.Range("N34").Value = -42876.668
$ws.Range("H76").Value = 42340.668
$ws.Range("J76").Value = 42340.668
$ws.Range("L76").Value = 42340.668
$ws.Range("N76").Value = -42970.668
$ws.Range("H79").Value = 42340.668
$ws.Range("J79").Value = 42340.668
$ws.Range("L79").Value = 42340.668
$ws.Range("N79").Value = -44524.668
$ws.Range("H80").Value = 2500.5715
$ws.Range("J80").Value = 2624.75
$ws.Range("L80").Value = 2624.75
$ws.Range("N80").Value = -4620.75
$ws.Range("H83").Value = 2500.5715
$ws.Range("J83").Value = 2624.75
$ws.Range("L83").Value = 13123.75
$ws.Range("N83").Value = -23107.75
$ws.Range("H126").Value = 14858.333
$ws.Range("I126").Value = 18066.928
$ws.Range("K126").Value = 54200.784
$ws.Range("M126").Value = -51730.784
$ws.Range("H132").Value = 2370.4285
$ws.Range("I132").Value = 2198.963
$ws.Range("K132").Value = 6596.889000000001
$ws.Range("M132").Value = -4066.889000000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4221.4287
$ws.Range("I132").Value = 4085.8096
$ws.Range("J132").Value = 4628.2856
$ws.Range("K132").Value = 12257.4288
$ws.Range("L132").Value = 13884.8568
$ws.Range("M132").Value = -9727.4288
$ws.Range("N132").Value = -18944.8568
$ws.Range("H136").Value = 3874.7368
$ws.Range("J136").Value = 4750
$ws.Range("L136").Value = 14250
$ws.Range("N136").Value = -19350

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H126").Value = 3436.75
$ws.Range("I126").Value = 3473.75
$ws.Range("K126").Value = 10421.25
$ws.Range("M126").Value = -7951.25
$ws.Range("H132").Value = 2298.1738
$ws.Range("I132").Value = 2380.7222
$ws.Range("K132").Value = 7142.1666
$ws.Range("M132").Value = -4612.1666
$ws.Range("H136").Value = 2755.5
$ws.Range("I136").Value = 2737.64
$ws.Range("K136").Value = 8212.92
$ws.Range("M136").Value = -5662.92
